$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: append 5 freshly-scraped rows at the bottom of the existing
# table (rows 70-74, before any insertion shifts things down). Copy the
# format of the last existing row first so date-looking text (e.g.
# "6/5/2024") is stored as literal text instead of being auto-converted
# into a date value.
$dates1 = @("5/52024", "6/5/2024", "7/5/2024", "8/5/2024", "9/5/2024")
$bvals1 = @(271891, 276455, 269692, 275778, 285859)
$cvals1 = @(215745, 215925, 213855, 209558, 209261)

for ($i = 0; $i -lt 5; $i++) {
    $r = 70 + $i
    $ws.Range("A69:C69").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $dates1[$i]
    $ws.Range("B$r").Value = $bvals1[$i]
    $ws.Range("C$r").Value = $cvals1[$i]
}

# --- Step 2: a missed day (26/4/2024) is found and inserted in its correct
# chronological slot, at row 62 - this shifts every row below it down by
# one (old rows 62-74 become 63-75). Insert() inherits the formatting of
# the row above, matching the target style.
$ws.Rows.Item(62).Insert()
$ws.Range("A62").Value = "26/4/2024"
$ws.Range("B62").Value = 296955
$ws.Range("C62").Value = 203915

# Stray empty, right-aligned style-only cells left over at G64/G65 (the
# old rows 63/64, now shifted to 64/65). Copy style from the adjacent A
# column cell (left-aligned Arial) and flip alignment to right.
$ws.Range("A64").Copy()
$ws.Range("G64").PasteSpecial(-4122)
$ws.Range("G64").HorizontalAlignment = -4152
$ws.Range("G64").Value = ""

$ws.Range("A65").Copy()
$ws.Range("G65").PasteSpecial(-4122)
$ws.Range("G65").HorizontalAlignment = -4152
$ws.Range("G65").Value = ""

# --- Step 3: one more day scraped afterwards, appended at the new bottom
# of the table (row 76).
$ws.Range("A75:C75").Copy()
$ws.Range("A76:C76").PasteSpecial(-4122)
$ws.Range("A76").Value = "10/5/2024"
$ws.Range("B76").Value = 297189
$ws.Range("C76").Value = 210616

# Update the selection to match (scroll position / topLeftCell is a
# view-only cosmetic attribute not exposed by this host).
$ws.Range("A76").Select()
